$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 271
$ws.Range("I33").Value = 276.7857
$ws.Range("K33").Value = 276.7857
$ws.Range("M33").Value = -47.78570000000002
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H64").Value = 8328.25
$ws.Range("I64").Value = 6609.75
$ws.Range("K64").Value = 6609.75
$ws.Range("M64").Value = -6361.75
$ws.Range("H67").Value = 8328.25
$ws.Range("I67").Value = 6609.75
$ws.Range("K67").Value = 6609.75
$ws.Range("M67").Value = -5751.75
$ws.Range("H98").Value = 1588259.1
$ws.Range("I98").Value = 1852818.4
$ws.Range("K98").Value = 1852818.4
$ws.Range("M98").Value = -1851320.4
$ws.Range("H120").Value = 48250
$ws.Range("J120").Value = 48250
$ws.Range("L120").Value = 48250
$ws.Range("N120").Value = -57926
$ws.Range("H122").Value = 1588259.1
$ws.Range("I122").Value = 1852818.4
$ws.Range("K122").Value = 5558455.199999999
$ws.Range("M122").Value = -5556005.199999999
$ws.Range("H132").Value = 1558.2307
$ws.Range("I132").Value = 1239.7391
$ws.Range("K132").Value = 3719.2173
$ws.Range("M132").Value = -1189.2173
$ws.Range("H137").Value = 3788.1428
$ws.Range("I137").Value = 3797.3333
$ws.Range("K137").Value = 11391.9999
$ws.Range("M137").Value = -8841.999899999999
$ws.Range("H138").Value = 6716.875
$ws.Range("I138").Value = 4873.5
$ws.Range("J138").Value = 7085.55
$ws.Range("K138").Value = 14620.5
$ws.Range("L138").Value = 21256.65
$ws.Range("M138").Value = -9480.5
$ws.Range("N138").Value = -31536.65

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19026014
$ws.Range("I32").Value = 19623618
$ws.Range("K32").Value = 19623618
$ws.Range("M32").Value = -19623331
$ws.Range("H132").Value = 4654.951
$ws.Range("I132").Value = 4111.25
$ws.Range("J132").Value = 5597.3667
$ws.Range("K132").Value = 12333.75
$ws.Range("L132").Value = 16792.1001
$ws.Range("M132").Value = -9803.75
$ws.Range("N132").Value = -21852.1001

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 34862.5
$ws.Range("J81").Value = 34862.5
$ws.Range("L81").Value = 34862.5
$ws.Range("N81").Value = -36984.5
$ws.Range("H84").Value = 34862.5
$ws.Range("J84").Value = 34862.5
$ws.Range("L84").Value = 104587.5
$ws.Range("N84").Value = -115195.5
$ws.Range("H86").Value = 5402.294
$ws.Range("I86").Value = 5449.1
$ws.Range("K86").Value = 5449.1
$ws.Range("M86").Value = -4326.1
$ws.Range("H89").Value = 5402.294
$ws.Range("I89").Value = 5449.1
$ws.Range("K89").Value = 27245.5
$ws.Range("M89").Value = -21629.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5975.7104
$ws.Range("I31").Value = 4290.154
$ws.Range("J31").Value = 6852.2
$ws.Range("K31").Value = 4290.154
$ws.Range("L31").Value = 6852.2
$ws.Range("M31").Value = -3995.154
$ws.Range("N31").Value = -7442.2
$ws.Range("H34").Value = 5975.7104
$ws.Range("I34").Value = 4290.154
$ws.Range("J34").Value = 6852.2
$ws.Range("K34").Value = 4290.154
$ws.Range("L34").Value = 6852.2
$ws.Range("M34").Value = -4088.154
$ws.Range("N34").Value = -7256.2
$ws.Range("H58").Value = 6202.8887
$ws.Range("I58").Value = 4456.75
$ws.Range("K58").Value = 4456.75
$ws.Range("M58").Value = -4253.75
$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25630
$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -27184
$ws.Range("H94").Value = 9232.143
$ws.Range("I94").Value = 18504.334
$ws.Range("J94").Value = 2278
$ws.Range("K94").Value = 18504.334
$ws.Range("L94").Value = 2278
$ws.Range("M94").Value = -18053.334
$ws.Range("N94").Value = -3180
$ws.Range("H132").Value = 3611
$ws.Range("I132").Value = 2900
$ws.Range("K132").Value = 8700
$ws.Range("M132").Value = -6170
$ws.Range("H134").Value = 4074.205
$ws.Range("I134").Value = 2419.3845
$ws.Range("K134").Value = 7258.1535
$ws.Range("M134").Value = -4723.1535
$ws.Range("H136").Value = 6202.8887
$ws.Range("I136").Value = 4456.75
$ws.Range("K136").Value = 13370.25
$ws.Range("M136").Value = -10820.25
$ws.Range("H141").Value = 73567
$ws.Range("J141").Value = 73567
$ws.Range("L141").Value = 73567
$ws.Range("N141").Value = -83927

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("M83").ClearContents()
$ws.Range("H107").Value = 552.375
$ws.Range("J107").Value = 730.4666999999999
$ws.Range("L107").Value = 2191.4001
$ws.Range("N107").Value = -6031.4001
$ws.Range("H114").Value = 1791.0769
$ws.Range("J114").Value = 2059.6
$ws.Range("L114").Value = 6178.799999999999
$ws.Range("N114").Value = -12686.8
$ws.Range("H122").Value = 3875
$ws.Range("J122").Value = 3500
$ws.Range("L122").Value = 31500
$ws.Range("N122").Value = -36400

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 27843980
$ws.Range("I80").Value = 106238.63
$ws.Range("J80").Value = 71431860
$ws.Range("K80").Value = 106238.63
$ws.Range("L80").Value = 71431860
$ws.Range("M80").Value = -105240.63
$ws.Range("N80").Value = -71433856
$ws.Range("H83").Value = 27843980
$ws.Range("I83").Value = 106238.63
$ws.Range("J83").Value = 71431860
$ws.Range("K83").Value = 531193.15
$ws.Range("L83").Value = 357159300
$ws.Range("M83").Value = -526201.15
$ws.Range("N83").Value = -357169284

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 13294.357
$ws.Range("I40").Value = 21051.834
$ws.Range("J40").Value = 7476.25
$ws.Range("K40").Value = 21051.834
$ws.Range("L40").Value = 7476.25
$ws.Range("M40").Value = -20915.834
$ws.Range("N40").Value = -7748.25
$ws.Range("H122").Value = 6984.2
$ws.Range("I122").Value = 4195.8887
$ws.Range("J122").Value = 11166.667
$ws.Range("K122").Value = 12587.6661
$ws.Range("L122").Value = 33500.001
$ws.Range("M122").Value = -10137.6661
$ws.Range("N122").Value = -38400.001
$ws.Range("H136").Value = 4504.1714
$ws.Range("I136").Value = 3578.4644
$ws.Range("K136").Value = 10735.3932
$ws.Range("M136").Value = -8185.393199999999

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4025
$ws.Range("H65").Value = 4025
$ws.Range("H126").Value = 2790.5557
$ws.Range("I126").Value = 2678.6667
$ws.Range("J126").Value = 3350
$ws.Range("K126").Value = 8036.000100000001
$ws.Range("L126").Value = 10050
$ws.Range("M126").Value = -5566.000100000001
$ws.Range("N126").Value = -14990
$ws.Range("H132").Value = 3776.6719
$ws.Range("I132").Value = 2767.2856
$ws.Range("J132").Value = 5703.6816
$ws.Range("K132").Value = 8301.856800000001
$ws.Range("L132").Value = 17111.0448
$ws.Range("M132").Value = -5771.856800000001
$ws.Range("N132").Value = -22171.0448
